$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text " Successful procedures were marked" " It is noticeably imbalanced with concern to target class; Successful procedures were marked"
Replace-Text " and unsuccessful" " (130 – 88.4%) and unsuccessful"
Replace-Text " as False." " as class False (17 – 11.6%)."
Replace-Text "Obtained results include " "Experiments included "
Replace-Text "100 % specificity" "100% specificity"
Replace-Text " score of 63 % on " " score of 63% on "
Replace-Text ". Accuracy was 87 %, while balanced (macro-average) accuracy was 92 %. P" ". Accuracy was 87%, while balanced accuracy was 92%. P"
Replace-Text " 46 %, but still the best" " 46%, but it was still the best compromise"
Replace-Text "extracted from d" "extracted from the d"
Replace-Text "as the most significant." "as the most significant ones."
Replace-Text "Other methods were also used, including " "Other methods of statistical analysis and machine learning were also used to identify important features, including "
Replace-Text "with their " "based on their "
